$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.676.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.363.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.34"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.50"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.379.46"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.35"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.120"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.428"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.966.68"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000185"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.69"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.823.12"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.389.15"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.17"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.87"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.36"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.91"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.993"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.10"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.522"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000113"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +16.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.29"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.79%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.96"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.96"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.40"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.32"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "22.87"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.65"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.88"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.43"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0759"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.886.19"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.80"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.79"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0313"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.84"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.25"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.740"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.81"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.05"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.08"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +12.93%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.28"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.822"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.38%  "
